$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (mac-address / document type test data rows)
$rows = @(
    @(10002, 10032, 3000176),
    @(10002, 10032, 3000177),
    @(10002, 10032, 3000178),
    @(10002, 10032, 3000179),
    @(10002, 10032, 3000180)
)

$startRow = 157
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
}

$ws.Range("E157").Select()
